# Updates Spriggan Profits workbook cell values across multiple sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market data,
# per the "chore: update Sheets via scheduled runner" commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4854.522
$ws.Range("J51").Value = 4692.6113
$ws.Range("L51").Value = 4692.6113
$ws.Range("N51").Value = -5660.6113
$ws.Range("H74").Value = 47709236
$ws.Range("I74").Value = 47709236
$ws.Range("K74").Value = 47709236
$ws.Range("M74").Value = -47708300
$ws.Range("H75").Value = 46437.668
$ws.Range("J75").Value = 46437.668
$ws.Range("L75").Value = 46437.668
$ws.Range("N75").Value = -48309.668
$ws.Range("H77").Value = 47709236
$ws.Range("I77").Value = 47709236
$ws.Range("K77").Value = 238546180
$ws.Range("M77").Value = -238541500
$ws.Range("H78").Value = 46437.668
$ws.Range("J78").Value = 46437.668
$ws.Range("L78").Value = 139313.004
$ws.Range("N78").Value = -148673.004
$ws.Range("H98").Value = 1746.8276
$ws.Range("I98").Value = 1584.9231
$ws.Range("K98").Value = 1584.9231
$ws.Range("M98").Value = -86.92309999999998
$ws.Range("H122").Value = 1746.8276
$ws.Range("I122").Value = 1584.9231
$ws.Range("K122").Value = 4754.7693
$ws.Range("M122").Value = -2304.7693
$ws.Range("H131").Value = 2325
$ws.Range("I131").Value = 1302.5714
$ws.Range("K131").Value = 3907.7142
$ws.Range("M131").Value = 1132.2858
$ws.Range("H137").Value = 2371.3513
$ws.Range("I137").Value = 2078.44
$ws.Range("J137").Value = 2981.5833
$ws.Range("K137").Value = 6235.32
$ws.Range("L137").Value = 8944.749899999999
$ws.Range("M137").Value = -3685.32
$ws.Range("N137").Value = -14044.7499
$ws.Range("H141").Value = 826.4
$ws.Range("I141").Value = 826.4
$ws.Range("K141").Value = 2479.2
$ws.Range("M141").Value = 2700.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2761.9312
$ws.Range("I32").Value = 2817.6086
$ws.Range("J32").Value = 2548.5
$ws.Range("K32").Value = 2817.6086
$ws.Range("L32").Value = 2548.5
$ws.Range("M32").Value = -2530.6086
$ws.Range("N32").Value = -3122.5
$ws.Range("H61").Value = 43479710
$ws.Range("I61").Value = 50001344
$ws.Range("J61").Value = 2166.3333
$ws.Range("K61").Value = 50001344
$ws.Range("L61").Value = 2166.3333
$ws.Range("M61").Value = -50001132
$ws.Range("N61").Value = -2590.3333
$ws.Range("H63").Value = 3089.2727
$ws.Range("I63").Value = 3089.2727
$ws.Range("K63").Value = 3089.2727
$ws.Range("M63").Value = -2403.2727
$ws.Range("H66").Value = 3089.2727
$ws.Range("I66").Value = 3089.2727
$ws.Range("K66").Value = 15446.3635
$ws.Range("M66").Value = -12014.3635
$ws.Range("H88").Value = 167567.5
$ws.Range("J88").Value = 1066.3334
$ws.Range("L88").Value = 1066.3334
$ws.Range("N88").Value = -1878.3334
$ws.Range("H91").Value = 167567.5
$ws.Range("J91").Value = 1066.3334
$ws.Range("L91").Value = 1066.3334
$ws.Range("N91").Value = -3874.3334
$ws.Range("H132").Value = 1888925.9
$ws.Range("I132").Value = 2224279.5
$ws.Range("K132").Value = 6672838.5
$ws.Range("M132").Value = -6670308.5
$ws.Range("H136").Value = 43479710
$ws.Range("I136").Value = 50001344
$ws.Range("J136").Value = 2166.3333
$ws.Range("K136").Value = 150004032
$ws.Range("L136").Value = 6498.999899999999
$ws.Range("M136").Value = -150001482
$ws.Range("N136").Value = -11598.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 72500
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 90000
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 90000
$ws.Range("M35").Value = -19690
$ws.Range("N35").Value = -90620
$ws.Range("H82").Value = 45999.4
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40766
$ws.Range("H85").Value = 45999.4
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8164.8647
$ws.Range("I31").Value = 6331.647
$ws.Range("J31").Value = 9723.1
$ws.Range("K31").Value = 6331.647
$ws.Range("L31").Value = 9723.1
$ws.Range("M31").Value = -6036.647
$ws.Range("N31").Value = -10313.1
$ws.Range("H34").Value = 8164.8647
$ws.Range("I34").Value = 6331.647
$ws.Range("J34").Value = 9723.1
$ws.Range("K34").Value = 6331.647
$ws.Range("L34").Value = 9723.1
$ws.Range("M34").Value = -6129.647
$ws.Range("N34").Value = -10127.1
$ws.Range("H132").Value = 21740850
$ws.Range("J132").Value = 1262.25
$ws.Range("L132").Value = 3786.75
$ws.Range("N132").Value = -8846.75
$ws.Range("H134").Value = 4168009.5
$ws.Range("I134").Value = 4718277
$ws.Range("K134").Value = 14154831
$ws.Range("M134").Value = -14152296

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 91571.63
$ws.Range("I5").Value = 143277
$ws.Range("J5").Value = 1087.25
$ws.Range("K5").Value = 429831
$ws.Range("L5").Value = 3261.75
$ws.Range("M5").Value = -429719
$ws.Range("N5").Value = -3485.75
$ws.Range("H32").Value = 1535.8
$ws.Range("J32").Value = 1535.8
$ws.Range("L32").Value = 4607.4
$ws.Range("N32").Value = -5173.4
$ws.Range("H33").Value = 2106.2856
$ws.Range("J33").Value = 2165.6667
$ws.Range("L33").Value = 12994.0002
$ws.Range("N33").Value = -13560.0002
$ws.Range("H68").Value = 4623.5
$ws.Range("J68").Value = 4747
$ws.Range("L68").Value = 14241
$ws.Range("N68").Value = -15863
$ws.Range("H71").Value = 4623.5
$ws.Range("J71").Value = 4747
$ws.Range("L71").Value = 42723
$ws.Range("N71").Value = -50835
$ws.Range("H80").Value = 9997.5
$ws.Range("I80").Value = 9997.5
$ws.Range("K80").Value = 29992.5
$ws.Range("M80").Value = -29056.5
$ws.Range("H83").Value = 9997.5
$ws.Range("I83").Value = 9997.5
$ws.Range("K83").Value = 89977.5
$ws.Range("M83").Value = -85297.5
$ws.Range("H135").Value = 91571.63
$ws.Range("I135").Value = 143277
$ws.Range("J135").Value = 1087.25
$ws.Range("K135").Value = 1289493
$ws.Range("L135").Value = 9785.25
$ws.Range("M135").Value = -1286958
$ws.Range("N135").Value = -14855.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 647.7
$ws.Range("J97").Value = 1599.3334
$ws.Range("L97").Value = 1599.3334
$ws.Range("N97").Value = -2591.3334
$ws.Range("H102").Value = 6722.2354
$ws.Range("I102").Value = 1618.5333
$ws.Range("K102").Value = 1618.5333
$ws.Range("M102").Value = 3.466699999999946
$ws.Range("H122").Value = 4262.525
$ws.Range("I122").Value = 2716.111
$ws.Range("K122").Value = 8148.333
$ws.Range("M122").Value = -5698.333
$ws.Range("H132").Value = 3788947.8
$ws.Range("I132").Value = 4033331.8
$ws.Range("K132").Value = 12099995.4
$ws.Range("M132").Value = -12097465.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6813.4546
$ws.Range("I40").Value = 6794.8
$ws.Range("K40").Value = 6794.8
$ws.Range("M40").Value = -6658.8
$ws.Range("H132").Value = 12912553
$ws.Range("I132").Value = 14203393
$ws.Range("K132").Value = 42610179
$ws.Range("M132").Value = -42607649

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 31839.4
$ws.Range("I2").Value = 32300
$ws.Range("K2").Value = 32300
$ws.Range("M2").Value = -32188
$ws.Range("H4").Value = 42882184
$ws.Range("I4").Value = 43750
$ws.Range("K4").Value = 43750
$ws.Range("M4").Value = -43637
$ws.Range("H122").Value = 1259.2858
$ws.Range("I122").Value = 1259.2858
$ws.Range("K122").Value = 3777.8574
$ws.Range("M122").Value = -1327.8574
$ws.Range("H132").Value = 11908772
$ws.Range("I132").Value = 14287217
$ws.Range("K132").Value = 42861651
$ws.Range("M132").Value = -42859121

